$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update cell text values (sharedStrings content) for the "semana del 6
#    hasta el 10 de marzo" block (rows 49-64). The two brand-new strings
#    ("parcial física mecánica" then "parcial cálculo integral") are written
#    first so they land in the shared-strings table in that order.
# ---------------------------------------------------------------------------
$ws.Range("E57").Value = "parcial física mecánica"
$ws.Range("F49").Value = "parcial cálculo integral"

$ws.Range("E49").Value = "clase de infórmatica"

$ws.Range("E50").Value = "clase de infórmatica"
$ws.Range("F50").Value = "clase de algebra lineal"

$ws.Range("E51").Value = "estudiar física mecánica"
$ws.Range("F51").Value = "parcial cálculo integral"

$ws.Range("E52").Value = "estudiar física mecánica"
$ws.Range("F52").Value = "parcial cálculo integral"

$ws.Range("F53").Value = "estudiar/prácticas/tareas informática"

$ws.Range("D55").Value = "estudiar algebra lineal"
$ws.Range("E55").Value = "estudiar física mecánica"

$ws.Range("D56").Value = "estudiar algebra lineal"
$ws.Range("E56").Value = "estudiar cálculo integral"
$ws.Range("F56").Value = "estudiar/prácticas/tareas informática"

$ws.Range("D57").Value = "Adelantar prácticas laboratorio informatica"

$ws.Range("E58").Value = "parcial física mecánica"

$ws.Range("D60").Value = "Adelantar prácticas laboratorio informatica"

$ws.Range("D61").Value = "Adelantar prácticas laboratorio informatica"

$ws.Range("D62").Value = "Adelantar prácticas laboratorio informatica"
$ws.Range("E62").Value = "estudiar cálculo integral"

$ws.Range("C63").Value = "estudiar física mecánica"
$ws.Range("D63").Value = "estudiar física mecánica"
$ws.Range("E63").Value = "estudiar cálculo integral"

$ws.Range("C64").Value = "estudiar física mecánica"
$ws.Range("E64").Value = "estudiar cálculo integral"

# ---------------------------------------------------------------------------
# 2) Fix up cell shading/format to match the new schedule. Rather than set
#    Interior colors directly (which would create brand-new fill/style
#    entries), copy the *format only* from a same-styled neighbour cell that
#    keeps its original formatting, so the existing style entries are reused.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# Cells that become "black" (busy/class block) shaded, bordered cells.
$blackDonor = $ws.Range("G49")
$blackTargets = @("E49","F49","E50","F50","F53","D55","D56","E56","C63","C64")
foreach ($ref in $blackTargets) {
    $blackDonor.Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
}

# Cells that become "white" shaded, bordered cells.
$whiteDonor = $ws.Range("F61")
$whiteTargets = @("D62","E62","D63")
foreach ($ref in $whiteTargets) {
    $whiteDonor.Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
}

# Cells that become plain (no fill) bordered cells.
$noneDonor = $ws.Range("C49")
$noneTargets = @("E51","F51","E52","F52","E55","F56","D57","E57","E58","D60","D61","E63","E64")
foreach ($ref in $noneTargets) {
    $noneDonor.Copy()
    $ws.Range($ref).PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Update the view state: scroll position + active selection.
# ---------------------------------------------------------------------------
$ws.Range("B51").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 53
$win.ScrollColumn = 1
